$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date style) of the last existing row (A374) down
# into the new rows A375:A385 before writing values, so the new date
# cells inherit the same style (s="2") as the rest of column A.
$ws.Range("A374").Copy($ws.Range("A375:A385"))

$data = @(
    @(44449, 0, 2, 47.65308553728854),
    @(44450, 1, 3, 71.47962830593281),
    @(44451, 0, 3, 71.47962830593281),
    @(44452, 1, 2, 47.65308553728854),
    @(44453, 1, 3, 71.47962830593281),
    @(44454, 0, 3, 71.47962830593281),
    @(44455, 0, 3, 71.47962830593281),
    @(44456, 0, 3, 71.47962830593281),
    @(44457, 0, 2, 47.65308553728854),
    @(44458, 0, 2, 47.65308553728854),
    @(44459, 0, 1, 23.82654276864427)
)

$row = 375
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
